$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 (Battery Cap): add a note in F11
$ws.Range("F11").Value = "*0.08"

# Row 12 (Male header Pins): clear quantity -> subtotal recalculates to 0
$ws.Range("B12").Value = ""

# Row 13 (Female Header Pins): clear quantity -> subtotal recalculates to 0
$ws.Range("B13").Value = ""

# Update selection to reflect last edited cell
$ws.Range("F12").Select()
